# Add homework rows for 2020-02-27 and 2020-02-28 (rows 63 and 64)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value forcing text storage (not auto-converted to
# date/number), then reset the style back to the default "Normal" style so
# the new cells don't pick up a stray number-format style.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 63: 2020-02-27
Set-TextCell 63 2 "2020-02-27"
Set-TextCell 63 3 "0215"
$ws.Cells.Item(63,1).Value = 1582761600
$ws.Cells.Item(63,4).Value = "SLVEST"
$ws.Cells.Item(63,5).Value = 1.31
$ws.Cells.Item(63,6).Value = 1.32
$ws.Cells.Item(63,7).Value = 1.17
$ws.Cells.Item(63,8).Value = 1.2
$ws.Cells.Item(63,9).Value = 19453400

# Row 64: 2020-02-28
Set-TextCell 64 2 "2020-02-28"
Set-TextCell 64 3 "0215"
$ws.Cells.Item(64,1).Value = 1582848000
$ws.Cells.Item(64,4).Value = "SLVEST"
$ws.Cells.Item(64,5).Value = 1.23
$ws.Cells.Item(64,6).Value = 1.25
$ws.Cells.Item(64,7).Value = 1.08
$ws.Cells.Item(64,8).Value = 1.14
$ws.Cells.Item(64,9).Value = 14910200
